$p = $ppt.ActivePresentation

# Slide 7 ("Proposed Solution") holds the two placeholders touched by this edit.
$s = $p.Slides.Item(7)

# Locate the shapes by name so the script is resilient to any shape-order
# differences, falling back to the known indices used in the source deck.
$titleShape = $null
$contentShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.Name -eq "Title 1") {
        $titleShape = $shape
    } elseif ($shape.Name -eq "Content Placeholder 2") {
        $contentShape = $shape
    }
}
if ($titleShape -eq $null) { $titleShape = $s.Shapes.Item(1) }
if ($contentShape -eq $null) { $contentShape = $s.Shapes.Item(2) }

# Title placeholder was an empty paragraph; it now reads "Proposed Solution".
$titleShape.TextFrame.TextRange.Text = "Proposed Solution"

# Extend Ryan's bullet with the additional mechanical/Fastenal detail.
$contentShape.TextFrame.TextRange.Text = "Ryan – systems (basic talk of what we want), unknowns, ways of testing systems, mention mechanical and what they think would work (Like Fastenal 3D array lockbox)"
